$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 25,14
$arr[0,0] = 61.06878433333333
$arr[0,1] = 183.206353
$arr[0,2] = 0.4308066250287063
$arr[0,3] = 0.4308066250287063
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 6.111751666666666
$arr[0,7] = 18.335255
$arr[0,8] = 0.6061514841909396
$arr[0,9] = 0.6061514841909394
$arr[0,10] = 373.2372444305572
$arr[0,11] = 3359.135199875015
$arr[0,12] = 0.2611340751604399
$arr[0,13] = 0.2611340751604398
$arr[1,0] = 61.06878433333333
$arr[1,1] = 183.206353
$arr[1,2] = 0.4308066250287063
$arr[1,3] = 0.4308066250287063
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 2.754304
$arr[1,7] = 8.262912
$arr[1,8] = 0.2731664420559804
$arr[1,9] = 0.2731664420559804
$arr[1,10] = 168.2019969644373
$arr[1,11] = 1513.817972679936
$arr[1,12] = 0.1176819129732366
$arr[1,13] = 0.1176819129732366
$arr[2,0] = 61.06878433333333
$arr[2,1] = 183.206353
$arr[2,2] = 0.4308066250287063
$arr[2,3] = 0.4308066250287063
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 0.568439
$arr[2,7] = 1.705317
$arr[2,8] = 0.0563766596410053
$arr[2,9] = 0.05637665964100529
$arr[2,10] = 34.71387869765567
$arr[2,11] = 312.424908278901
$arr[2,12] = 0.02428743847033357
$arr[2,13] = 0.02428743847033357
$arr[3,0] = 61.06878433333333
$arr[3,1] = 183.206353
$arr[3,2] = 0.4308066250287063
$arr[3,3] = 0.4308066250287063
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 0.3689163333333333
$arr[3,7] = 1.106749
$arr[3,8] = 0.036588394815171
$arr[3,9] = 0.036588394815171
$arr[3,10] = 22.52927199737745
$arr[3,11] = 202.763447976397
$arr[3,12] = 0.01576252288554164
$arr[3,13] = 0.01576252288554164
$arr[4,0] = 61.06878433333333
$arr[4,1] = 183.206353
$arr[4,2] = 0.4308066250287063
$arr[4,3] = 0.4308066250287063
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.2794673333333333
$arr[4,7] = 0.838402
$arr[4,8] = 0.02771701929690382
$arr[4,9] = 0.02771701929690381
$arr[4,10] = 17.06673030754511
$arr[4,11] = 153.600572767906
$arr[4,12] = 0.01194067553915466
$arr[4,13] = 0.01194067553915466
$arr[5,0] = 3.621603
$arr[5,1] = 10.864809
$arr[5,2] = 0.02554841368886107
$arr[5,3] = 0.02554841368886107
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 6.111751666666666
$arr[5,7] = 18.335255
$arr[5,8] = 0.6061514841909396
$arr[5,9] = 0.6061514841909394
$arr[5,10] = 22.134338171255
$arr[5,11] = 199.209043541295
$arr[5,12] = 0.01548620887622726
$arr[5,13] = 0.01548620887622725
$arr[6,0] = 3.621603
$arr[6,1] = 10.864809
$arr[6,2] = 0.02554841368886107
$arr[6,3] = 0.02554841368886107
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 2.754304
$arr[6,7] = 8.262912
$arr[6,8] = 0.2731664420559804
$arr[6,9] = 0.2731664420559804
$arr[6,10] = 9.974995629312
$arr[6,11] = 89.774960663808
$arr[6,12] = 0.006978969267560485
$arr[6,13] = 0.006978969267560484
$arr[7,0] = 3.621603
$arr[7,1] = 10.864809
$arr[7,2] = 0.02554841368886107
$arr[7,3] = 0.02554841368886107
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 0.568439
$arr[7,7] = 1.705317
$arr[7,8] = 0.0563766596410053
$arr[7,9] = 0.05637665964100529
$arr[7,10] = 2.058660387717
$arr[7,11] = 18.527943489453
$arr[7,12] = 0.001440334222904522
$arr[7,13] = 0.001440334222904521
$arr[8,0] = 3.621603
$arr[8,1] = 10.864809
$arr[8,2] = 0.02554841368886107
$arr[8,3] = 0.02554841368886107
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 0.3689163333333333
$arr[8,7] = 1.106749
$arr[8,8] = 0.036588394815171
$arr[8,9] = 0.036588394815171
$arr[8,10] = 1.336068499549
$arr[8,11] = 12.024616495941
$arr[8,12] = 0.0009347754469493685
$arr[8,13] = 0.0009347754469493684
$arr[9,0] = 3.621603
$arr[9,1] = 10.864809
$arr[9,2] = 0.02554841368886107
$arr[9,3] = 0.02554841368886107
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 0.2794673333333333
$arr[9,7] = 0.838402
$arr[9,8] = 0.02771701929690382
$arr[9,9] = 0.02771701929690381
$arr[9,10] = 1.012119732802
$arr[9,11] = 9.109077595218
$arr[9,12] = 0.000708125875219444
$arr[9,13] = 0.0007081258752194439
$arr[10,0] = 36.843258
$arr[10,1] = 110.529774
$arr[10,2] = 0.2599088848306786
$arr[10,3] = 0.2599088848306786
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 6.111751666666666
$arr[10,7] = 18.335255
$arr[10,8] = 0.6061514841909396
$arr[10,9] = 0.6061514841909394
$arr[10,10] = 225.17684348693
$arr[10,11] = 2026.59159138237
$arr[10,12] = 0.1575441562945278
$arr[10,13] = 0.1575441562945278
$arr[11,0] = 36.843258
$arr[11,1] = 110.529774
$arr[11,2] = 0.2599088848306786
$arr[11,3] = 0.2599088848306786
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 2.754304
$arr[11,7] = 8.262912
$arr[11,8] = 0.2731664420559804
$arr[11,9] = 0.2731664420559804
$arr[11,10] = 101.477532882432
$arr[11,11] = 913.297795941888
$arr[11,12] = 0.07099838532793407
$arr[11,13] = 0.07099838532793407
$arr[12,0] = 36.843258
$arr[12,1] = 110.529774
$arr[12,2] = 0.2599088848306786
$arr[12,3] = 0.2599088848306786
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 0.568439
$arr[12,7] = 1.705317
$arr[12,8] = 0.0563766596410053
$arr[12,9] = 0.05637665964100529
$arr[12,10] = 20.943144734262
$arr[12,11] = 188.488302608358
$arr[12,12] = 0.01465279473777241
$arr[12,13] = 0.01465279473777241
$arr[13,0] = 36.843258
$arr[13,1] = 110.529774
$arr[13,2] = 0.2599088848306786
$arr[13,3] = 0.2599088848306786
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 0.3689163333333333
$arr[13,7] = 1.106749
$arr[13,8] = 0.036588394815171
$arr[13,9] = 0.036588394815171
$arr[13,10] = 13.592079649414
$arr[13,11] = 122.328716844726
$arr[13,12] = 0.009509648894155681
$arr[13,13] = 0.009509648894155681
$arr[14,0] = 36.843258
$arr[14,1] = 110.529774
$arr[14,2] = 0.2599088848306786
$arr[14,3] = 0.2599088848306786
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 0.2794673333333333
$arr[14,7] = 0.838402
$arr[14,8] = 0.02771701929690382
$arr[14,9] = 0.02771701929690381
$arr[14,10] = 10.296487064572
$arr[14,11] = 92.66838358114801
$arr[14,12] = 0.007203899576288672
$arr[14,13] = 0.007203899576288671
$arr[15,0] = 2.119603
$arr[15,1] = 6.358808999999999
$arr[15,2] = 0.01495263127961596
$arr[15,3] = 0.01495263127961596
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 6.111751666666666
$arr[15,7] = 18.335255
$arr[15,8] = 0.6061514841909396
$arr[15,9] = 0.6061514841909394
$arr[15,10] = 12.95448716792166
$arr[15,11] = 116.590384511295
$arr[15,12] = 0.00906355964269908
$arr[15,13] = 0.009063559642699079
$arr[16,0] = 2.119603
$arr[16,1] = 6.358808999999999
$arr[16,2] = 0.01495263127961596
$arr[16,3] = 0.01495263127961596
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 2.754304
$arr[16,7] = 8.262912
$arr[16,8] = 0.2731664420559804
$arr[16,9] = 0.2731664420559804
$arr[16,10] = 5.838031021311999
$arr[16,11] = 52.54227919180799
$arr[16,12] = 0.004084557086027652
$arr[16,13] = 0.004084557086027652
$arr[17,0] = 2.119603
$arr[17,1] = 6.358808999999999
$arr[17,2] = 0.01495263127961596
$arr[17,3] = 0.01495263127961596
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 0.568439
$arr[17,7] = 1.705317
$arr[17,8] = 0.0563766596410053
$arr[17,9] = 0.05637665964100529
$arr[17,10] = 1.204865009717
$arr[17,11] = 10.843785087453
$arr[17,12] = 0.0008429794043883583
$arr[17,13] = 0.0008429794043883582
$arr[18,0] = 2.119603
$arr[18,1] = 6.358808999999999
$arr[18,2] = 0.01495263127961596
$arr[18,3] = 0.01495263127961596
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 0.3689163333333333
$arr[18,7] = 1.106749
$arr[18,8] = 0.036588394815171
$arr[18,9] = 0.036588394815171
$arr[18,10] = 0.7819561668823333
$arr[18,11] = 7.037605501940999
$arr[18,12] = 0.0005470927767842642
$arr[18,13] = 0.0005470927767842642
$arr[19,0] = 2.119603
$arr[19,1] = 6.358808999999999
$arr[19,2] = 0.01495263127961596
$arr[19,3] = 0.01495263127961596
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 0.2794673333333333
$arr[19,7] = 0.838402
$arr[19,8] = 0.02771701929690382
$arr[19,9] = 0.02771701929690381
$arr[19,10] = 0.5923597981353332
$arr[19,11] = 5.331238183217999
$arr[19,12] = 0.0004144423697166031
$arr[19,13] = 0.000414442369716603
$arr[20,0] = 38.101267
$arr[20,1] = 114.303801
$arr[20,2] = 0.268783445172138
$arr[20,3] = 0.268783445172138
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 6.111751666666666
$arr[20,7] = 18.335255
$arr[20,8] = 0.6061514841909396
$arr[20,9] = 0.6061514841909394
$arr[20,10] = 232.8654820893616
$arr[20,11] = 2095.789338804254
$arr[20,12] = 0.1629234842170455
$arr[20,13] = 0.1629234842170454
$arr[21,0] = 38.101267
$arr[21,1] = 114.303801
$arr[21,2] = 0.268783445172138
$arr[21,3] = 0.268783445172138
$arr[21,4] = 3
$arr[21,5] = 1
$arr[21,6] = 2.754304
$arr[21,7] = 8.262912
$arr[21,8] = 0.2731664420559804
$arr[21,9] = 0.2731664420559804
$arr[21,10] = 104.942472103168
$arr[21,11] = 944.4822489285119
$arr[21,12] = 0.07342261740122163
$arr[21,13] = 0.07342261740122161
$arr[22,0] = 38.101267
$arr[22,1] = 114.303801
$arr[22,2] = 0.268783445172138
$arr[22,3] = 0.268783445172138
$arr[22,4] = 3
$arr[22,5] = 1
$arr[22,6] = 0.568439
$arr[22,7] = 1.705317
$arr[22,8] = 0.0563766596410053
$arr[22,9] = 0.05637665964100529
$arr[22,10] = 21.658246112213
$arr[22,11] = 194.924215009917
$arr[22,12] = 0.01515311280560644
$arr[22,13] = 0.01515311280560643
$arr[23,0] = 38.101267
$arr[23,1] = 114.303801
$arr[23,2] = 0.268783445172138
$arr[23,3] = 0.268783445172138
$arr[23,4] = 3
$arr[23,5] = 1
$arr[23,6] = 0.3689163333333333
$arr[23,7] = 1.106749
$arr[23,8] = 0.036588394815171
$arr[23,9] = 0.036588394815171
$arr[23,10] = 14.05617971699433
$arr[23,11] = 126.505617452949
$arr[23,12] = 0.009834354811740055
$arr[23,13] = 0.009834354811740054
$arr[24,0] = 38.101267
$arr[24,1] = 114.303801
$arr[24,2] = 0.268783445172138
$arr[24,3] = 0.268783445172138
$arr[24,4] = 3
$arr[24,5] = 1
$arr[24,6] = 0.2794673333333333
$arr[24,7] = 0.838402
$arr[24,8] = 0.02771701929690382
$arr[24,9] = 0.02771701929690381
$arr[24,10] = 10.64805948511133
$arr[24,11] = 95.83253536600199
$arr[24,12] = 0.01194067553915466
$arr[24,13] = 0.01194067553915466
$ws.Range("G2:T26").Value = $arr
Write-Output "done"
